# Daily attendance processing - 2026-01-15 15:08:00
# Swap the order of "<user>, System" -> "System, <user>" in column G
# for the specific rows touched by this processing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows whose column G value uses "admin@admin.com, System"
$adminRows = @(7, 33, 59)
foreach ($r in $adminRows) {
    $ws.Cells.Item($r, 7).Value = "System, admin@admin.com"
}

# Rows whose column G value uses "dnasr281@gmail.com, System"
$dnasrRows = @(
    10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26,
    36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52,
    62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78,
    83, 84, 85, 86, 90, 92, 93, 94, 96, 99, 101,
    109, 110, 111, 112, 116, 118, 119, 120, 122, 125, 127,
    135, 136, 137, 138, 142, 144, 145, 146, 148, 151, 153
)
foreach ($r in $dnasrRows) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}
